$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "25.838.45"

Set-TextValue $ws.Range("D3") "1.630.50"
Set-TextValue $ws.Range("E3") "  +0.09%  "

Set-TextValue $ws.Range("D5") "214.13"
Set-TextValue $ws.Range("E5") "  +0.19%  "

Set-TextValue $ws.Range("E6") "  +0.99%  "

Set-TextValue $ws.Range("E7") "  +0.59%  "

Set-TextValue $ws.Range("E8") "  -0.46%  "

Set-TextValue $ws.Range("E9") "  +0.21%  "

Set-TextValue $ws.Range("D10") "19.53"
Set-TextValue $ws.Range("E10") "  -0.44%  "

Set-TextValue $ws.Range("E11") "  +0.31%  "

Set-TextValue $ws.Range("D12") "1.856.08"
Set-TextValue $ws.Range("E12") "  +0.11%  "

Set-TextValue $ws.Range("D13") "4.25"
Set-TextValue $ws.Range("E13") "  -0.11%  "

Set-TextValue $ws.Range("D14") "1.622.89"
Set-TextValue $ws.Range("E14") "  -0.60%  "

Set-TextValue $ws.Range("E15") "  -1.37%  "

Set-TextValue $ws.Range("E16") "  -0.42%  "

Set-TextValue $ws.Range("D17") "62.57"
Set-TextValue $ws.Range("E17") "  -0.05%  "

Set-TextValue $ws.Range("D18") "25.845.05"
Set-TextValue $ws.Range("E18") "  +0.20%  "

Set-TextValue $ws.Range("E19") "  +0.58%  "

Set-TextValue $ws.Range("B20") "Uniswap"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D20") "4.39"
Set-TextValue $ws.Range("E20") "  -1.01%  "

Set-TextValue $ws.Range("B21") "BitcoinCash"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D21") "193.10"
Set-TextValue $ws.Range("E21") "  +1.28%  "

Set-TextValue $ws.Range("D22") "9.91"
Set-TextValue $ws.Range("E22") "  -0.12%  "

Set-TextValue $ws.Range("E23") "  -0.48%  "

Set-TextValue $ws.Range("E24") "  +1.18%  "

Set-TextValue $ws.Range("D25") "143.26"
Set-TextValue $ws.Range("E25") "  +0.74%  "

Set-TextValue $ws.Range("E26") "  +0.59%  "

Set-TextValue $ws.Range("E27") "  +2.69%  "

Set-TextValue $ws.Range("D28") "6.83"
Set-TextValue $ws.Range("E28") "  -0.04%  "

Set-TextValue $ws.Range("D29") "15.40"
Set-TextValue $ws.Range("E29") "  -0.62%  "

Set-TextValue $ws.Range("E30") "  +0.28%  "

Set-TextValue $ws.Range("D31") "0.0497"
Set-TextValue $ws.Range("E31") "  +0.74%  "

Set-TextValue $ws.Range("E32") "  -0.76%  "

Set-TextValue $ws.Range("E33") "  -0.02%  "

Set-TextValue $ws.Range("E34") "  -2.04%  "

Set-TextValue $ws.Range("E35") "  +1.69%  "

Set-TextValue $ws.Range("E36") "  -0.14%  "

Set-TextValue $ws.Range("D37") "1.138.53"
Set-TextValue $ws.Range("E37") "  -0.13%  "

Set-TextValue $ws.Range("E38") "  +0.43%  "

Set-TextValue $ws.Range("E39") "  -1.00%  "

Set-TextValue $ws.Range("E40") "  +0.49%  "

Set-TextValue $ws.Range("E41") "  +0.70%  "

Set-TextValue $ws.Range("D42") "99.09"
Set-TextValue $ws.Range("E42") "  -1.23%  "

Set-TextValue $ws.Range("E43") "  -2.93%  "

Set-TextValue $ws.Range("D44") "0.796"
Set-TextValue $ws.Range("E44") "  -0.41%  "

Set-TextValue $ws.Range("D45") "1.765.74"

Set-TextValue $ws.Range("B46") "Aave"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D46") "56.17"
Set-TextValue $ws.Range("E46") "  +1.45%  "

Set-TextValue $ws.Range("B47") "Cronos"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D47") "0.0528"
Set-TextValue $ws.Range("E47") "  +3.12%  "

Set-TextValue $ws.Range("B48") "Mantle"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D48") "0.415"
Set-TextValue $ws.Range("E48") "  -0.12%  "

Set-TextValue $ws.Range("B49") "RenderToken"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D49") "1.43"
Set-TextValue $ws.Range("E49") "  -2.08%  "

Set-TextValue $ws.Range("B50") "EnergySwap"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "7.63"
Set-TextValue $ws.Range("E50") "  +0.83%  "

Set-TextValue $ws.Range("B51") "Algorand"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.0958"
Set-TextValue $ws.Range("E51") "  +0.23%  "
